$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price (D) and Volume(1h) (E) values for each data row (2-51).
$updates = @(
    @{ Row = 2; D = '23.910.49'; E = '  +16.28%  ' }
    @{ Row = 3; D = '1.659.36'; E = '  +12.86%  ' }
    @{ Row = 4; D = '0.9945'; E = '  -1.57%  ' }
    @{ Row = 5; D = '307.30'; E = '  +11.06%  ' }
    @{ Row = 6; D = '0.9898'; E = '  +3.23%  ' }
    @{ Row = 7; D = '0.3722'; E = '  +4.61%  ' }
    @{ Row = 8; D = '0.3450'; E = '  +12.80%  ' }
    @{ Row = 9; D = '47.46'; E = '  +21.22%  ' }
    @{ Row = 10; D = '1.175'; E = '  +8.61%  ' }
    @{ Row = 11; D = '0.07205'; E = '  +8.97%  ' }
    @{ Row = 12; D = '0.9907'; E = '  -1.31%  ' }
    @{ Row = 13; D = '20.64'; E = '  +14.51%  ' }
    @{ Row = 14; D = '6.000'; E = '  +10.21%  ' }
    @{ Row = 15; D = '6.721'; E = '  +9.05%  ' }
    @{ Row = 16; D = '1.656.27'; E = '  +12.88%  ' }
    @{ Row = 17; D = '0.00001096'; E = '  +7.82%  ' }
    @{ Row = 18; D = '0.9884'; E = '  +2.97%  ' }
    @{ Row = 19; D = '0.06749'; E = '  +13.78%  ' }
    @{ Row = 20; D = '80.94'; E = '  +17.86%  ' }
    @{ Row = 21; D = '16.41'; E = '  +13.75%  ' }
    @{ Row = 22; D = '6.072'; E = '  +11.19%  ' }
    @{ Row = 23; D = '11.93'; E = '  +6.97%  ' }
    @{ Row = 24; D = '23.832.17'; E = '  +15.85%  ' }
    @{ Row = 25; D = '2.345'; E = '  +3.05%  ' }
    @{ Row = 26; D = '3.395'; E = '  -8.76%  ' }
    @{ Row = 27; D = '2.673'; E = '  +28.64%  ' }
    @{ Row = 28; D = '151.65'; E = '  +4.06%  ' }
    @{ Row = 29; D = '19.53'; E = '  +14.50%  ' }
    @{ Row = 30; D = '1.838.77'; E = '  +12.94%  ' }
    @{ Row = 31; D = '126.34'; E = '  +10.64%  ' }
    @{ Row = 32; D = '4.053'; E = '  +1.33%  ' }
    @{ Row = 33; D = '6.143'; E = '  +25.31%  ' }
    @{ Row = 34; D = '0.9799'; E = '  +24.45%  ' }
    @{ Row = 35; D = '1.706'; E = '  +17.47%  ' }
    @{ Row = 36; D = '0.08367'; E = '  +5.60%  ' }
    @{ Row = 37; D = '12.24'; E = '  +19.62%  ' }
    @{ Row = 38; D = '8.880'; E = '  +22.24%  ' }
    @{ Row = 39; D = '0.06332'; E = '  +11.84%  ' }
    @{ Row = 40; D = '5.275'; E = '  +11.90%  ' }
    @{ Row = 41; D = '1.279'; E = '  +5.51%  ' }
    @{ Row = 42; D = '0.02289'; E = '  +13.13%  ' }
    @{ Row = 43; D = '0.2063'; E = '  +12.07%  ' }
    @{ Row = 44; D = '0.6049'; E = '  +16.15%  ' }
    @{ Row = 45; D = '0.9856'; E = '  +2.62%  ' }
    @{ Row = 46; D = '3.821'; E = '  +8.83%  ' }
    @{ Row = 47; D = '13.11'; E = '  +8.86%  ' }
    @{ Row = 48; D = '0.5898'; E = '  +14.45%  ' }
    @{ Row = 49; D = '127.04'; E = '  +6.30%  ' }
    @{ Row = 50; D = '1.995'; E = '  +11.00%  ' }
    @{ Row = 51; D = '0.07035'; E = '  +9.62%  ' }
)

foreach ($u in $updates) {
    # Prefix D values with a leading apostrophe so Excel stores them as
    # literal text (matching the original inline-string cells) instead of
    # auto-converting numeric-looking strings (e.g. "6.000", "1.175") into numbers.
    $dCell = $ws.Cells.Item($u.Row, 4)
    $dCell.Value = "'" + $u.D
    $dCell.Style = "Normal"

    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
